$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 2.1937118031827127
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.9851839106752087
$ws.Range("E2").ClearContents()

# Row 3 values
$ws.Range("B3").Value = 1.944872166344664
$ws.Range("C3").Value = -0.90133221548553433
$ws.Range("D3").Value = 2.8889582240326792
$ws.Range("E3").Value = -0.79127357551109223

# Update the selected range to reflect the reduced data block (B1:E3)
$ws.Range("B1:E3").Select()
